# Update the cryptos list sheet with the latest scraped values.
# Commit: Updated cryptos list on Wed Feb 14 15:26:43 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "Price" values in column D are stored as plain text in the workbook
# (e.g. "51.695.57" or "0.0₃0962"), not real numbers. Some of the new values
# look like ordinary decimals (e.g. "117.09") and Excel would silently store
# those as numbers instead of text unless the cell is explicitly marked as
# Text first. Cells whose new text isn't number-like (contains more than one
# "." or special glyphs) are untouched, since Excel already keeps those as
# text without any extra nudging.
$numericLookingPriceCells = @(
  "D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D17","D19",
  "D20","D21","D23","D24","D26","D29","D31","D32","D33","D34","D35","D36",
  "D37","D38","D39","D40","D41","D42","D43","D45","D46","D48","D50","D51"
)
foreach ($cellRef in $numericLookingPriceCells) {
  $ws.Range($cellRef).NumberFormat = "@"
}

# --- Rows whose B (Coin) / C (Link) / D (Price) / E (Volume 1h) values change ---

# Row 11: Chainlink -> Dogecoin
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.0833"
$ws.Range("E11").Value = "  +3.04%  "

# Row 12: Dogecoin -> Chainlink
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "19.96"
$ws.Range("E12").Value = "  +1.20%  "

# Row 42: EnergySwap -> VeChain
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0347"
$ws.Range("E42").Value = "  +11.35%  "

# Row 43: VeChain -> EnergySwap
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "23.16"
$ws.Range("E43").Value = "  +2.58%  "

# --- Rows whose D (Price) and E (Volume 1h) values change ---

$ws.Range("D2").Value = "51.695.57"
$ws.Range("E2").Value = "  +6.08%  "

$ws.Range("D3").Value = "2.751.36"
$ws.Range("E3").Value = "  +4.81%  "

$ws.Range("D5").Value = "117.09"
$ws.Range("E5").Value = "  +6.30%  "

$ws.Range("D6").Value = "331.52"
$ws.Range("E6").Value = "  +2.90%  "

$ws.Range("D7").Value = "0.533"
$ws.Range("E7").Value = "  +2.52%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "0.576"
$ws.Range("E9").Value = "  +7.00%  "

$ws.Range("D10").Value = "41.53"
$ws.Range("E10").Value = "  +5.40%  "

$ws.Range("D13").Value = "0.129"
$ws.Range("E13").Value = "  +2.69%  "

$ws.Range("D14").Value = "7.62"
$ws.Range("E14").Value = "  +5.92%  "

$ws.Range("D15").Value = "3.175.29"
$ws.Range("E15").Value = "  +5.12%  "

$ws.Range("D16").Value = "2.757.59"
$ws.Range("E16").Value = "  +5.08%  "

$ws.Range("D17").Value = "0.883"
$ws.Range("E17").Value = "  +2.96%  "

$ws.Range("D18").Value = "51.556.89"
$ws.Range("E18").Value = "  +5.98%  "

$ws.Range("D19").Value = "13.55"
$ws.Range("E19").Value = "  +5.72%  "

$ws.Range("D20").Value = "3.03"
$ws.Range("E20").Value = "  +4.44%  "

$ws.Range("D21").Value = "6.82"
$ws.Range("E21").Value = "  +2.48%  "

$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  +2.41%  "

$ws.Range("D23").Value = "278.93"
$ws.Range("E23").Value = "  +3.09%  "

$ws.Range("D24").Value = "69.65"
$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("D26").Value = "26.82"
$ws.Range("E26").Value = "  +3.10%  "

$ws.Range("D29").Value = "10.23"
$ws.Range("E29").Value = "  +1.63%  "

$ws.Range("D31").Value = "0.140"
$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("D32").Value = "34.93"
$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").Value = "50.33"
$ws.Range("E33").Value = "  +1.86%  "

$ws.Range("D34").Value = "5.58"
$ws.Range("E34").Value = "  +3.20%  "

$ws.Range("D35").Value = "0.0821"
$ws.Range("E35").Value = "  +3.59%  "

$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").Value = "19.07"
$ws.Range("E37").Value = "  -0.52%  "

$ws.Range("D38").Value = "2.09"
$ws.Range("E38").Value = "  +3.38%  "

$ws.Range("D39").Value = "4.94"
$ws.Range("E39").Value = "  +0.59%  "

$ws.Range("D40").Value = "3.20"
$ws.Range("E40").Value = "  +2.40%  "

$ws.Range("D41").Value = "130.73"
$ws.Range("E41").Value = "  +5.15%  "

$ws.Range("D45").Value = "2.27"
$ws.Range("E45").Value = "  +6.54%  "

$ws.Range("D46").Value = "2.39"
$ws.Range("E46").Value = "  +14.63%  "

$ws.Range("D47").Value = "2.107.86"
$ws.Range("E47").Value = "  +1.99%  "

$ws.Range("D48").Value = "3.33"
$ws.Range("E48").Value = "  +4.14%  "

$ws.Range("D50").Value = "5.56"
$ws.Range("E50").Value = "  +8.08%  "

$ws.Range("D51").Value = "8.95"
$ws.Range("E51").Value = "  +0.44%  "

# --- Row 4 changes only D (Price); E4 stays the same ---
$ws.Range("D4").Value = "0.999"

# --- Rows whose only E (Volume 1h) changes; D stays the same ---
$ws.Range("E25").Value = "  +4.54%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E49").Value = "  +3.01%  "
